$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.818.66"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.309.75"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.68"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.38"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.76"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.24"
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.80"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "2.664.25"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "2.310.67"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("D18").Value = "42.769.61"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.26"
$ws.Range("E19").Value = "  -5.29%  "
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.85"
$ws.Range("E23").Value = "  +8.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.58"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.48"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  +15.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.81"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.08"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.23"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.58"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  -6.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0703"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.75"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.28"
$ws.Range("E42").Value = "  +11.31%  "
$ws.Range("D43").Value = "1.976.44"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.45"
$ws.Range("E44").Value = "  +5.01%  "
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "2.534.75"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.46"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.79"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.58"
$ws.Range("E51").Value = "  +0.44%  "
